# Updated via Streamlit Approval System
#
# - Resets the COST_CENTER / LEDGER_NAME / LEDGER_UNDER / TO / BY testing
#   placeholder values (columns AK:AO) back to numeric 0 for every data row.
# - Clears the stray APPROVAL_1 / APPROVAL_2 (columns AI:AJ) test markers
#   (HOLD / PAID / ACCEPTED) back to blank, keeping them as empty text
#   cells (not fully-blank) to match the sheet's normal "empty string"
#   convention for untouched approval cells.
# - Corrects BASIC_AMOUNT on row 12 from 7001 back to 7000.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 27

for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Range("AK$row").Value = 0
    $ws.Range("AL$row").Value = 0
    $ws.Range("AM$row").Value = 0
    $ws.Range("AN$row").Value = 0
    $ws.Range("AO$row").Value = 0
}

# Cells that held leftover approval test text - blank them out to an empty
# (text-typed) string, the same representation the sheet already uses for
# every other un-set approval cell.
$approvalCellsToClear = @(
    "AI3",
    "AI4", "AJ4",
    "AI5", "AJ5",
    "AI7",
    "AI9", "AJ9",
    "AI10",
    "AI12", "AJ12",
    "AI22",
    "AI23",
    "AI27"
)

foreach ($cellRef in $approvalCellsToClear) {
    $ws.Range($cellRef).Value = "'"
    $ws.Range($cellRef).Style = "Normal"
}

# Fix the BASIC_AMOUNT typo on row 12 (7001 -> 7000).
$ws.Range("V12").Value = 7000
